$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 19351.666
$ws.Range("I38").Value = 36.666668
$ws.Range("K38").Value = 110.000004
$ws.Range("M38").Value = 261.999996
$ws.Range("H40").Value = 128740.75
$ws.Range("I40").Value = 1501200
$ws.Range("K40").Value = 1501200
$ws.Range("M40").Value = -1501025
$ws.Range("H101").Value = 282
$ws.Range("I101").Value = 265.14285
$ws.Range("K101").Value = 795.4285500000001
$ws.Range("M101").Value = 826.5714499999999
$ws.Range("H133").Value = 124998.5
$ws.Range("J133").Value = 124998.5
$ws.Range("L133").Value = 124998.5
$ws.Range("N133").Value = -135118.5
$ws.Range("H137").Value = 4359.7095
$ws.Range("I137").Value = 3803.25
$ws.Range("J137").Value = 5371.4546
$ws.Range("K137").Value = 11409.75
$ws.Range("L137").Value = 16114.3638
$ws.Range("M137").Value = -8859.75
$ws.Range("N137").Value = -21214.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6124.75
$ws.Range("I45").Value = 5166.5
$ws.Range("J45").Value = 8999.5
$ws.Range("K45").Value = 5166.5
$ws.Range("L45").Value = 8999.5
$ws.Range("M45").Value = -4789.5
$ws.Range("N45").Value = -9753.5
$ws.Range("H61").Value = 1909.4348
$ws.Range("I61").Value = 1772.4186
$ws.Range("J61").Value = 3873.3333
$ws.Range("K61").Value = 1772.4186
$ws.Range("L61").Value = 3873.3333
$ws.Range("M61").Value = -1560.4186
$ws.Range("N61").Value = -4297.3333
$ws.Range("H103").Value = 78329.5
$ws.Range("J103").Value = 78329.5
$ws.Range("L103").Value = 78329.5
$ws.Range("N103").Value = -80673.5
$ws.Range("H110").Value = 1860.8
$ws.Range("I110").Value = 1321.8
$ws.Range("J110").Value = 2399.8
$ws.Range("K110").Value = 1321.8
$ws.Range("L110").Value = 2399.8
$ws.Range("M110").Value = 723.2
$ws.Range("N110").Value = -6489.8
$ws.Range("H122").Value = 2502.5151
$ws.Range("I122").Value = 2147.8708
$ws.Range("K122").Value = 6443.6124
$ws.Range("M122").Value = -3993.6124
$ws.Range("H136").Value = 1909.4348
$ws.Range("I136").Value = 1772.4186
$ws.Range("J136").Value = 3873.3333
$ws.Range("K136").Value = 5317.2558
$ws.Range("L136").Value = 11619.9999
$ws.Range("M136").Value = -2767.2558
$ws.Range("N136").Value = -16719.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 198.42857
$ws.Range("I22").Value = 207.5
$ws.Range("J22").Value = 144
$ws.Range("K22").Value = 207.5
$ws.Range("L22").Value = 144
$ws.Range("M22").Value = -34.5
$ws.Range("N22").Value = -490
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 509999.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 509999.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -512121.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 509999.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 1529998.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -1540606.5
$ws.Range("H94").Value = 2033.6666
$ws.Range("I94").Value = 829.7143
$ws.Range("J94").Value = 6247.5
$ws.Range("K94").Value = 829.7143
$ws.Range("L94").Value = 6247.5
$ws.Range("M94").Value = -378.7143
$ws.Range("N94").Value = -7149.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4325.48
$ws.Range("I31").Value = 2697.9092
$ws.Range("J31").Value = 5604.2856
$ws.Range("K31").Value = 2697.9092
$ws.Range("L31").Value = 5604.2856
$ws.Range("M31").Value = -2402.9092
$ws.Range("N31").Value = -6194.2856
$ws.Range("H34").Value = 4325.48
$ws.Range("I34").Value = 2697.9092
$ws.Range("J34").Value = 5604.2856
$ws.Range("K34").Value = 2697.9092
$ws.Range("L34").Value = 5604.2856
$ws.Range("M34").Value = -2495.9092
$ws.Range("N34").Value = -6008.2856
$ws.Range("H75").Value = 104643.336
$ws.Range("J75").Value = 104643.336
$ws.Range("L75").Value = 104643.336
$ws.Range("N75").Value = -106639.336
$ws.Range("H78").Value = 104643.336
$ws.Range("J78").Value = 104643.336
$ws.Range("L78").Value = 313930.008
$ws.Range("N78").Value = -323914.008
$ws.Range("H99").Value = 1989.5
$ws.Range("I99").Value = 1992.6666
$ws.Range("K99").Value = 1992.6666
$ws.Range("M99").Value = -494.6666
$ws.Range("H100").Value = 110790
$ws.Range("J100").Value = 110790
$ws.Range("L100").Value = 110790
$ws.Range("N100").Value = -112954
$ws.Range("H107").Value = 145358.72
$ws.Range("I107").Value = 251502.75
$ws.Range("K107").Value = 251502.75
$ws.Range("M107").Value = -249582.75
$ws.Range("H122").Value = 3103.1667
$ws.Range("I122").Value = 2695.2222
$ws.Range("K122").Value = 8085.6666
$ws.Range("M122").Value = -5635.6666
$ws.Range("H126").Value = 1989.5
$ws.Range("I126").Value = 1992.6666
$ws.Range("K126").Value = 5977.9998
$ws.Range("M126").Value = -3507.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.42856999999999
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 152.5
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 915
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -1141
$ws.Range("H4").Value = 75499770
$ws.Range("I4").Value = 46975670
$ws.Range("J4").Value = 225251250
$ws.Range("K4").Value = 140927010
$ws.Range("L4").Value = 675753750
$ws.Range("M4").Value = -140926898
$ws.Range("N4").Value = -675753974
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -9588
$ws.Range("H131").Value = 1712.8
$ws.Range("J131").Value = 1854.4286
$ws.Range("L131").Value = 5563.2858
$ws.Range("N131").Value = -15643.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2557
$ws.Range("J80").Value = 2766.5
$ws.Range("L80").Value = 2766.5
$ws.Range("N80").Value = -4762.5
$ws.Range("H83").Value = 2557
$ws.Range("J83").Value = 2766.5
$ws.Range("L83").Value = 13832.5
$ws.Range("N83").Value = -23816.5
$ws.Range("H128").Value = 134759
$ws.Range("J128").Value = 134759
$ws.Range("L128").Value = 134759
$ws.Range("N128").Value = -144719
$ws.Range("H132").Value = 3347.3044
$ws.Range("I132").Value = 3363.0908
$ws.Range("K132").Value = 10089.2724
$ws.Range("M132").Value = -7559.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19749.75
$ws.Range("I7").Value = 19499.5
$ws.Range("K7").Value = 19499.5
$ws.Range("M7").Value = -19387.5
$ws.Range("H22").Value = 2272.6667
$ws.Range("I22").Value = 2224.375
$ws.Range("K22").Value = 2224.375
$ws.Range("M22").Value = -1929.375
$ws.Range("H27").Value = 2272.6667
$ws.Range("I27").Value = 2224.375
$ws.Range("K27").Value = 2224.375
$ws.Range("M27").Value = -2117.375
$ws.Range("H46").Value = 5100.5
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 5400.7144
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 5400.7144
$ws.Range("M46").Value = -2811
$ws.Range("N46").Value = -5776.7144
$ws.Range("H106").Value = 1035565.2
$ws.Range("J106").Value = 1035565.2
$ws.Range("L106").Value = 1035565.2
$ws.Range("N106").Value = -1038089.2
$ws.Range("H126").Value = 19749.75
$ws.Range("I126").Value = 19499.5
$ws.Range("K126").Value = 58498.5
$ws.Range("M126").Value = -56028.5
$ws.Range("H136").Value = 3166.838
$ws.Range("I136").Value = 2380.5
$ws.Range("J136").Value = 5613.222
$ws.Range("K136").Value = 7141.5
$ws.Range("L136").Value = 16839.666
$ws.Range("M136").Value = -4591.5
$ws.Range("N136").Value = -21939.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 109000
$ws.Range("J75").Value = 109000
$ws.Range("L75").Value = 109000
$ws.Range("N75").Value = -110872
$ws.Range("H78").Value = 109000
$ws.Range("J78").Value = 109000
$ws.Range("L78").Value = 327000
$ws.Range("M78").Value = -336360
$ws.Range("H107").Value = 432.7
$ws.Range("I107").Value = 306.16666
$ws.Range("K107").Value = 918.4999799999999
$ws.Range("M107").Value = 1001.50002
$ws.Range("H126").Value = 11601.857
$ws.Range("I126").Value = 11868.833
$ws.Range("K126").Value = 35606.499
$ws.Range("M126").Value = -33136.499
$ws.Range("H130").Value = 88877
$ws.Range("J130").Value = 88877
$ws.Range("L130").Value = 88877
$ws.Range("N130").Value = -98917
$ws.Range("H136").Value = 14968.59
$ws.Range("I136").Value = 1337.1321
$ws.Range("J136").Value = 43867.28
$ws.Range("K136").Value = 4011.3963
$ws.Range("L136").Value = 131601.84
$ws.Range("M136").Value = -1461.3963
$ws.Range("N136").Value = -136701.84
